# "break out stock.yaml completed" - append the next breakout scan batch
# (24/06/2024 05:45:27) to the "10per change" sheet, and fix the two
# bsecode cells in the prior batch that were stored as text instead of
# numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

function Set-TextValue($cell, [string]$text) {
    # Force the cell to stay a text cell even when its content looks
    # numeric (matches the sheet's existing "bsecode" column convention).
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# E30 / E31 ("bsecode") were inline text; the edit normalizes them to
# real numbers, matching every other bsecode cell in the sheet.
$ws.Cells.Item(30, 5).Value = 20
$ws.Cells.Item(31, 5).Value = 531344

# New batch: row 32 (BSE), row 33 (CONCOR), row 34 (ZEEL)
$ws.Cells.Item(32, 1).Value = "24/06/2024 05:45:27"
$ws.Cells.Item(32, 2).Value = 1
$ws.Cells.Item(32, 3).Value = "BSE"
$ws.Cells.Item(32, 4).Value = "BSE (Bombay stock exchange)"
Set-TextValue $ws.Cells.Item(32, 5) "20"
$ws.Cells.Item(32, 6).Value = -2.06
$ws.Cells.Item(32, 7).Value = 2507.35
$ws.Cells.Item(32, 8).Value = 443628

$ws.Cells.Item(33, 1).Value = "24/06/2024 05:45:27"
$ws.Cells.Item(33, 2).Value = 2
$ws.Cells.Item(33, 3).Value = "CONCOR"
$ws.Cells.Item(33, 4).Value = "Container Corporation Of India Limited"
Set-TextValue $ws.Cells.Item(33, 5) "531344"
$ws.Cells.Item(33, 6).Value = -3.11
$ws.Cells.Item(33, 7).Value = 1056.95
$ws.Cells.Item(33, 8).Value = 1850473

$ws.Cells.Item(34, 1).Value = "24/06/2024 05:45:27"
$ws.Cells.Item(34, 2).Value = 3
$ws.Cells.Item(34, 3).Value = "ZEEL"
$ws.Cells.Item(34, 4).Value = "Zee Entertainment Enterprises Limited"
Set-TextValue $ws.Cells.Item(34, 5) "505537"
$ws.Cells.Item(34, 6).Value = -1.75
$ws.Cells.Item(34, 7).Value = 151.54
$ws.Cells.Item(34, 8).Value = 5174436
